$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A95").Value = "2025/12/06 18:00"
$ws.Range("B95").Value = "-"
$ws.Range("C95").Value = "-"
$ws.Range("D95").Value = "-"
$ws.Range("E95").Value = "-"
$ws.Range("F95").Value = "-"
$ws.Range("G95").Value = "-"
